# Fall 2022 schedule - week 7 (and week 8 for Safety Dance) results entered.
# Previously-unplayed weeks were marked "A" (not yet Available/Assigned) as a
# placeholder; now that the games happened, fill in the real outcomes
# (W = win, L = loss, DNP = did not play, NA = not applicable).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")

# --- "Wookie Mistakes" table (rows 3-10): week of 10/11 (column H) completed ---
$ws.Range("H3").Value  = "DNP"   # Daniel Burcham
$ws.Range("H4").Value  = "L"     # Leo Hayward
$ws.Range("H5").Value  = "L"     # Laura Thompson
$ws.Range("H6").Value  = "DNP"   # Kim Quan
$ws.Range("H7").Value  = "W"     # Scott Berry
$ws.Range("H8").Value  = "W"     # Jason Liess
$ws.Range("H9").Value  = "DNP"   # Jason Bohrer
$ws.Range("H10").Value = "L"     # Dan Aquino

# --- "Safety Dance" table (rows 15-22): weeks of 10/11 (H) and 10/18 (I) completed ---
$ws.Range("H15").Value = "W"     # Jason Bohrer
$ws.Range("I15").Value = "W"

$ws.Range("H16").Value = "DNP"   # Jason Liess
$ws.Range("I16").Value = "L"

$ws.Range("H17").Value = "W"     # Daniel Burcham
$ws.Range("I17").Value = "NA"

$ws.Range("H18").Value = "DNP"   # Scott Berry
$ws.Range("I18").Value = "NA"

$ws.Range("H19").Value = "W"     # Dan Aquino
$ws.Range("I19").Value = "NA"

$ws.Range("H20").Value = "DNP"   # Ashley Daniels
$ws.Range("I20").Value = "W"

$ws.Range("H21").Value = "L"     # Adrian Warden
$ws.Range("I21").Value = "W"

$ws.Range("H22").Value = "L"     # Shelia Lowe
$ws.Range("I22").Value = "L"

# Leave the cursor where the author's last edit landed.
$ws.Range("J17").Select()
